$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.Value = "'" + $text
    $r.Style = 'Normal'
}

Set-TextCell 'D2' '67.665.68'
Set-TextCell 'E2' '  -0.95%  '
Set-TextCell 'D3' '3.786.77'
Set-TextCell 'E3' '  +1.13%  '
Set-TextCell 'E4' '  -0.06%  '
Set-TextCell 'D5' '595.85'
Set-TextCell 'E5' '  +0.54%  '
Set-TextCell 'D6' '167.02'
Set-TextCell 'E6' '  +0.62%  '
Set-TextCell 'D7' '3.773.10'
Set-TextCell 'E7' '  +0.83%  '
Set-TextCell 'E8' '  +0.05%  '
Set-TextCell 'E9' '  +0.09%  '
Set-TextCell 'E10' '  -0.09%  '
Set-TextCell 'D11' '6.31'
Set-TextCell 'E11' '  -1.86%  '
Set-TextCell 'E12' '  +0.23%  '
Set-TextCell 'E13' '  -2.63%  '
Set-TextCell 'D14' '35.93'
Set-TextCell 'E14' '  -0.30%  '
Set-TextCell 'D15' '4.422.80'
Set-TextCell 'E15' '  +1.06%  '
Set-TextCell 'D16' '3.793.88'
Set-TextCell 'E16' '  +1.05%  '
Set-TextCell 'D17' '18.49'
Set-TextCell 'E17' '  +3.58%  '
Set-TextCell 'D18' '67.643.47'
Set-TextCell 'E18' '  -1.00%  '
Set-TextCell 'E19' '  +0.77%  '
Set-TextCell 'E20' '  -0.14%  '
Set-TextCell 'D21' '10.04'
Set-TextCell 'E21' '  -6.09%  '
Set-TextCell 'D22' '458.59'
Set-TextCell 'E22' '  -1.17%  '
Set-TextCell 'E23' '  +0.36%  '
Set-TextCell 'D24' '0.0000152'
Set-TextCell 'E24' '  +4.03%  '
Set-TextCell 'D25' '83.39'
Set-TextCell 'E25' '  -0.57%  '
Set-TextCell 'E26' '  +1.23%  '
Set-TextCell 'D27' '2.12'
Set-TextCell 'E27' '  -2.43%  '
Set-TextCell 'B28' 'Dai'
Set-TextCell 'C28' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 'D28' '1.00'
Set-TextCell 'E28' '  +0.14%  '
Set-TextCell 'B29' 'RenderToken'
Set-TextCell 'C29' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D29' '10.01'
Set-TextCell 'E29' '  -0.46%  '
Set-TextCell 'D30' '3.933.20'
Set-TextCell 'E30' '  +1.01%  '
Set-TextCell 'D31' '2.78'
Set-TextCell 'E31' '  +0.66%  '
Set-TextCell 'D32' '2.21'
Set-TextCell 'E32' '  +3.12%  '
Set-TextCell 'D33' '7.18'
Set-TextCell 'E33' '  -1.47%  '
Set-TextCell 'D34' '29.61'
Set-TextCell 'E34' '  -0.75%  '
Set-TextCell 'D35' '0.999'
Set-TextCell 'D36' '9.07'
Set-TextCell 'E36' '  -1.09%  '
Set-TextCell 'E37' '  -0.28%  '
Set-TextCell 'D38' '3.35'
Set-TextCell 'E38' '  -2.22%  '
Set-TextCell 'E39' '  +0.10%  '
Set-TextCell 'D40' '0.995'
Set-TextCell 'E40' '  -0.43%  '
Set-TextCell 'D41' '5.77'
Set-TextCell 'E41' '  +0.05%  '
Set-TextCell 'E42' '  -0.18%  '
Set-TextCell 'B43' 'USDe'
Set-TextCell 'C43' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextCell 'D43' '1.00'
Set-TextCell 'E43' '  +0.00%  '
Set-TextCell 'B44' 'Arweave'
Set-TextCell 'C44' 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextCell 'D44' '45.76'
Set-TextCell 'E44' '  +5.67%  '
Set-TextCell 'D45' '48.13'
Set-TextCell 'E45' '  +3.41%  '
Set-TextCell 'E46' '  -0.68%  '
Set-TextCell 'D47' '149.38'
Set-TextCell 'E48' '  -1.80%  '
Set-TextCell 'D49' '392.89'
Set-TextCell 'E49' '  +0.89%  '
Set-TextCell 'D50' '26.70'
Set-TextCell 'E50' '  +4.74%  '
Set-TextCell 'E51' '  -5.06%  '
